# Update "想去人数" (want-to-go count) figures in column F across all four sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 256
$ws1.Range("F5").Value  = 41
$ws1.Range("F6").Value  = 2134
$ws1.Range("F7").Value  = 222
$ws1.Range("F8").Value  = 662
$ws1.Range("F9").Value  = 27
$ws1.Range("F10").Value = 182
$ws1.Range("F11").Value = 150
$ws1.Range("F12").Value = 656
$ws1.Range("F13").Value = 48
$ws1.Range("F14").Value = 88
$ws1.Range("F15").Value = 1282
$ws1.Range("F19").Value = 246

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 20
$ws2.Range("F6").Value  = 10
$ws2.Range("F8").Value  = 8
$ws2.Range("F9").Value  = 118
$ws2.Range("F11").Value = 30
$ws2.Range("F12").Value = 209

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6322
$ws3.Range("F3").Value = 786
$ws3.Range("F4").Value = 1979
$ws3.Range("F5").Value = 210

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6322
$ws4.Range("F3").Value  = 786
$ws4.Range("F4").Value  = 1979
$ws4.Range("F6").Value  = 210
$ws4.Range("F10").Value = 20
$ws4.Range("F12").Value = 256
$ws4.Range("F13").Value = 44
$ws4.Range("F14").Value = 10
$ws4.Range("F16").Value = 8
$ws4.Range("F17").Value = 2134
$ws4.Range("F18").Value = 118
$ws4.Range("F19").Value = 222
$ws4.Range("F21").Value = 30
$ws4.Range("F22").Value = 662
$ws4.Range("F23").Value = 27
$ws4.Range("F24").Value = 182
$ws4.Range("F25").Value = 209
$ws4.Range("F26").Value = 150
$ws4.Range("F27").Value = 656
$ws4.Range("F28").Value = 48
$ws4.Range("F29").Value = 88
$ws4.Range("F31").Value = 1282
$ws4.Range("F43").Value = 246
